$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H131").Value = 2329.6667
$ws.Range("I131").Value = 2329.6667
$ws.Range("K131").Value = 6989.000100000001
$ws.Range("M131").Value = -1949.000100000001

$ws.Range("H137").Value = 1330.5
$ws.Range("I137").Value = 1497
$ws.Range("K137").Value = 4491
$ws.Range("M137").Value = -1941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3121.7932
$ws.Range("I32").Value = 2674.5
$ws.Range("J32").Value = 6998.3335
$ws.Range("K32").Value = 2674.5
$ws.Range("L32").Value = 6998.3335
$ws.Range("M32").Value = -2387.5
$ws.Range("N32").Value = -7572.3335

$ws.Range("H61").Value = 1024.5
$ws.Range("I61").Value = 1024.5
$ws.Range("K61").Value = 1024.5
$ws.Range("M61").Value = -812.5

$ws.Range("H74").Value = 1058.7778
$ws.Range("I74").Value = 1058.7778
$ws.Range("K74").Value = 1058.7778
$ws.Range("M74").Value = -184.7778000000001

$ws.Range("H77").Value = 1058.7778
$ws.Range("I77").Value = 1058.7778
$ws.Range("K77").Value = 5293.889
$ws.Range("M77").Value = -925.8890000000001

$ws.Range("H122").Value = 547.5
$ws.Range("I122").Value = 547.5
$ws.Range("K122").Value = 1642.5
$ws.Range("M122").Value = 807.5

$ws.Range("H132").Value = 2109.5557
$ws.Range("I132").Value = 2310.75
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 6932.25
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -4402.25
$ws.Range("N132").Value = -6560

$ws.Range("H136").Value = 1024.5
$ws.Range("I136").Value = 1024.5
$ws.Range("K136").Value = 3073.5
$ws.Range("M136").Value = -523.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 175.66667
$ws.Range("I107").Value = 175.66667
$ws.Range("K107").Value = 175.66667
$ws.Range("M107").Value = 1744.33333

$ws.Range("H134").Value = 3467.8635
$ws.Range("I134").Value = 3490.3809
$ws.Range("K134").Value = 10471.1427
$ws.Range("M134").Value = -7936.1427

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 874.75
$ws.Range("I58").Value = 666.3333
$ws.Range("K58").Value = 666.3333
$ws.Range("M58").Value = -463.3333

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H132").Value = 1453.6364
$ws.Range("I132").Value = 1453.6364
$ws.Range("K132").Value = 4360.9092
$ws.Range("M132").Value = -1830.9092

$ws.Range("H136").Value = 874.75
$ws.Range("I136").Value = 666.3333
$ws.Range("K136").Value = 1998.9999
$ws.Range("M136").Value = 551.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.5
$ws.Range("J2").Value = 142.83333
$ws.Range("L2").Value = 856.9999799999999
$ws.Range("N2").Value = -1082.99998

$ws.Range("H3").Value = 4499
$ws.Range("I3").Value = 4499
$ws.Range("K3").Value = 13497
$ws.Range("M3").Value = -13385

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9500
$ws.Range("I57").Value = 6000
$ws.Range("K57").Value = 6000
$ws.Range("M57").Value = -5180

$ws.Range("H80").Value = 25000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 25000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 6589.6665
$ws.Range("I126").Value = 6400
$ws.Range("J126").Value = 6969
$ws.Range("K126").Value = 19200
$ws.Range("L126").Value = 20907
$ws.Range("M126").Value = -16730
$ws.Range("N126").Value = -25847

$ws.Range("H132").Value = 4426.3887
$ws.Range("I132").Value = 3885.0667
$ws.Range("J132").Value = 7133
$ws.Range("K132").Value = 11655.2001
$ws.Range("L132").Value = 21399
$ws.Range("M132").Value = -9125.2001
$ws.Range("N132").Value = -26459

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2249.8333
$ws.Range("I40").Value = 2373.75
$ws.Range("K40").Value = 2373.75
$ws.Range("M40").Value = -2237.75

$ws.Range("H55").Value = 839.8
$ws.Range("I55").Value = 783.3333
$ws.Range("J55").Value = 924.5
$ws.Range("K55").Value = 783.3333
$ws.Range("L55").Value = 924.5
$ws.Range("M55").Value = -610.3333
$ws.Range("N55").Value = -1270.5

$ws.Range("H82").Value = 700
$ws.Range("I82").Value = 700
$ws.Range("K82").Value = 700
$ws.Range("M82").Value = -339

$ws.Range("H85").Value = 700
$ws.Range("I85").Value = 700
$ws.Range("K85").Value = 700
$ws.Range("M85").Value = 548

$ws.Range("H132").Value = 4007.3572
$ws.Range("I132").Value = 2665.4443
$ws.Range("K132").Value = 7996.3329
$ws.Range("M132").Value = -5466.3329

$ws.Range("H136").Value = 8712.857
$ws.Range("I136").Value = 8831.666999999999
$ws.Range("K136").Value = 26495.001
$ws.Range("M136").Value = -23945.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 4647900.5
$ws.Range("I100").Value = 5362137
$ws.Range("K100").Value = 10724274
$ws.Range("M100").Value = -10723733

$ws.Range("H107").Value = 1007.7
$ws.Range("I107").Value = 743.0769
$ws.Range("K107").Value = 2229.2307
$ws.Range("M107").Value = -309.2307000000001

$ws.Range("H122").Value = 1957
$ws.Range("I122").Value = 1957
$ws.Range("K122").Value = 5871
$ws.Range("M122").Value = -3421

$ws.Range("H132").Value = 2600.111
$ws.Range("I132").Value = 2434
$ws.Range("K132").Value = 7302
$ws.Range("M132").Value = -4772

$ws.Range("H136").Value = 1304.2222
$ws.Range("I136").Value = 1217.25
$ws.Range("K136").Value = 3651.75
$ws.Range("M136").Value = -1101.75
